$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the price cells we are about to edit to Text format so that
# values such as "174.10" or "66.001.81" are preserved exactly as
# literal text instead of being auto-converted to numbers.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

$ws.Range("D2").Value = "66.001.81"
$ws.Range("D3").Value = "3.445.32"
$ws.Range("E3").Value = "  -0.17%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "585.63"
$ws.Range("E5").Value = "  +0.88%  "
$ws.Range("D6").Value = "174.10"
$ws.Range("E6").Value = "  -1.00%  "
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("E8").Value = "  +0.51%  "
$ws.Range("D9").Value = "3.444.72"
$ws.Range("E9").Value = "  -0.17%  "
$ws.Range("E10").Value = "  -1.44%  "
$ws.Range("D11").Value = "6.95"
$ws.Range("E11").Value = "  +1.55%  "
$ws.Range("D12").Value = "0.418"
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").Value = "4.044.44"
$ws.Range("E14").Value = "  +1.80%  "
$ws.Range("D15").Value = "29.17"
$ws.Range("E15").Value = "  -4.35%  "
$ws.Range("D16").Value = "65.991.10"
$ws.Range("E16").Value = "  -0.80%  "
$ws.Range("E17").Value = "  -0.04%  "
$ws.Range("D18").Value = "3.445.32"
$ws.Range("E18").Value = "  -0.17%  "
$ws.Range("E19").Value = "  -0.74%  "
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D21").Value = "370.69"
$ws.Range("E21").Value = "  -1.38%  "
$ws.Range("E22").Value = "  -1.25%  "
$ws.Range("D23").Value = "72.50"
$ws.Range("E23").Value = "  +1.89%  "
$ws.Range("D24").Value = "0.999"
$ws.Range("E24").Value = "  +0.15%  "
$ws.Range("D25").Value = "0.533"
$ws.Range("E25").Value = "  +1.11%  "
$ws.Range("E26").Value = "  +4.13%  "
$ws.Range("E27").Value = "  -0.64%  "
$ws.Range("E28").Value = "  +3.89%  "
$ws.Range("E29").Value = "  -0.17%  "
$ws.Range("E30").Value = "  -0.72%  "
$ws.Range("E31").Value = "  +0.38%  "
$ws.Range("D32").Value = "23.63"
$ws.Range("E32").Value = "  -1.36%  "
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("D34").Value = "7.03"
$ws.Range("E34").Value = "  +0.00%  "
$ws.Range("E35").Value = "  -4.72%  "
$ws.Range("E36").Value = "  +1.25%  "
$ws.Range("D37").Value = "161.55"
$ws.Range("E37").Value = "  +1.33%  "
$ws.Range("D38").Value = "0.880"
$ws.Range("E38").Value = "  +0.38%  "
$ws.Range("D39").Value = "28.30"
$ws.Range("E39").Value = "  +3.42%  "
$ws.Range("E40").Value = "  +0.21%  "
$ws.Range("D41").Value = "2.60"
$ws.Range("E41").Value = "  -0.92%  "
$ws.Range("D42").Value = "2.790.06"
$ws.Range("E42").Value = "  +3.84%  "
$ws.Range("E43").Value = "  +0.12%  "
$ws.Range("E44").Value = "  -0.33%  "
$ws.Range("E45").Value = "  -0.96%  "
$ws.Range("D46").Value = "25.16"
$ws.Range("E46").Value = "  -0.22%  "
$ws.Range("D47").Value = "39.79"
$ws.Range("E47").Value = "  -0.94%  "
$ws.Range("E48").Value = "  -0.56%  "
$ws.Range("D49").Value = "328.36"
$ws.Range("E49").Value = "  +2.54%  "
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("E51").Value = "  +1.58%  "
